# Applies: "Add object detector selection feature"
# - Normalizes D263:D300 from text-encoded numbers to true numeric cells
# - Appends vehicle-speed log rows 301:370 captured across several detector runs
#   (rows 301:356 logged as numeric speeds, rows 357:370 logged as text speeds,
#   matching the source export behaviour being reproduced here)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: convert D263:D300 from text to numeric values ---
$dFix = [ordered]@{
    263 = "31.51"
    264 = "45.77"
    265 = "41.46"
    266 = "41.01"
    267 = "51.41"
    268 = "41.89"
    269 = "42.43"
    270 = "0"
    271 = "0"
    272 = "0"
    273 = "0"
    274 = "0"
    275 = "0"
    276 = "0"
    277 = "0"
    278 = "0"
    279 = "0"
    280 = "0"
    281 = "0"
    282 = "0"
    283 = "38.01"
    284 = "43.2"
    285 = "48.93"
    286 = "90.28"
    287 = "15.26"
    288 = "44.99"
    289 = "42.89"
    290 = "40.32"
    291 = "44.76"
    292 = "35.97"
    293 = "41.35"
    294 = "51.75"
    295 = "123.76"
    296 = "87.57"
    297 = "74.44"
    298 = "112.09"
    299 = "111.01"
    300 = "19.27"
}
foreach ($r in $dFix.Keys) {
    $ws.Cells.Item($r, 4).Value = [double]$dFix[$r]
}

# --- Part 2: append new rows 301:370 (object detector selection log rows) ---
# Columns: Row, Time, ID, D-cell-type ("n"=numeric, "s"=text), D-value (as text)
$newRows = @(
    ,@(301, "2025-06-03 22:25:32", 2, "n", "31.51")
    ,@(302, "2025-06-03 22:25:32", 4, "n", "45.77")
    ,@(303, "2025-06-03 22:25:33", 5, "n", "41.46")
    ,@(304, "2025-06-03 22:25:34", 8, "n", "18.28")
    ,@(305, "2025-06-03 22:25:38", 11, "n", "15.8")
    ,@(306, "2025-06-03 22:25:39", 9, "n", "41.01")
    ,@(307, "2025-06-03 22:25:39", 12, "n", "18.09")
    ,@(308, "2025-06-03 22:25:41", 10, "n", "51.41")
    ,@(309, "2025-06-03 22:25:42", 13, "n", "41.89")
    ,@(310, "2025-06-03 22:25:45", 14, "n", "42.43")
    ,@(311, "2025-06-03 22:25:45", 16, "n", "17.84")
    ,@(312, "2025-06-03 22:25:47", 17, "n", "17.59")
    ,@(313, "2025-06-03 22:32:03", 2, "n", "31.51")
    ,@(314, "2025-06-03 22:32:04", 4, "n", "45.77")
    ,@(315, "2025-06-03 22:32:05", 5, "n", "41.46")
    ,@(316, "2025-06-03 22:32:19", 12, "n", "0")
    ,@(317, "2025-06-03 22:32:20", 14, "n", "0")
    ,@(318, "2025-06-03 22:32:21", 15, "n", "0")
    ,@(319, "2025-06-03 22:32:22", 18, "n", "0")
    ,@(320, "2025-06-03 23:31:58", 2, "n", "31.51")
    ,@(321, "2025-06-03 23:31:59", 4, "n", "45.77")
    ,@(322, "2025-06-03 23:32:00", 5, "n", "41.46")
    ,@(323, "2025-06-03 23:32:01", 8, "n", "18.28")
    ,@(324, "2025-06-03 23:32:12", 11, "n", "27.01")
    ,@(325, "2025-06-03 23:32:15", 12, "n", "33.04")
    ,@(326, "2025-06-03 23:32:19", 17, "n", "40.92")
    ,@(327, "2025-06-03 23:32:25", 14, "n", "18")
    ,@(328, "2025-06-03 23:32:38", 30, "n", "17.16")
    ,@(329, "2025-06-03 23:32:39", 36, "n", "0")
    ,@(330, "2025-06-03 23:32:44", 19, "n", "41.03")
    ,@(331, "2025-06-03 23:32:47", 33, "n", "18.09")
    ,@(332, "2025-06-03 23:32:54", 20, "n", "58.54")
    ,@(333, "2025-06-03 23:32:59", 26, "n", "41.8")
    ,@(334, "2025-06-03 23:33:09", 32, "n", "42.34")
    ,@(335, "2025-06-03 23:33:09", 37, "n", "18.01")
    ,@(336, "2025-06-03 23:33:19", 39, "n", "17.52")
    ,@(337, "2025-06-03 23:33:24", 35, "n", "34.66")
    ,@(338, "2025-06-03 23:33:35", 41, "n", "19.11")
    ,@(339, "2025-06-03 23:33:45", 38, "n", "35.5")
    ,@(340, "2025-06-03 23:33:45", 43, "n", "16.67")
    ,@(341, "2025-06-03 23:33:52", 45, "n", "16.91")
    ,@(342, "2025-06-03 23:34:03", 40, "n", "34.7")
    ,@(343, "2025-06-03 23:34:23", 47, "n", "17.73")
    ,@(344, "2025-06-03 23:34:30", 48, "n", "16.66")
    ,@(345, "2025-06-03 23:34:37", 51, "n", "18.1")
    ,@(346, "2025-06-03 23:34:56", 49, "n", "39.92")
    ,@(347, "2025-06-03 23:35:04", 55, "n", "17.32")
    ,@(348, "2025-06-03 23:35:09", 56, "n", "19.23")
    ,@(349, "2025-06-03 23:35:11", 57, "n", "18.45")
    ,@(350, "2025-06-03 23:35:17", 52, "n", "38.61")
    ,@(351, "2025-06-03 23:35:40", 61, "n", "17.17")
    ,@(352, "2025-06-03 23:35:40", 62, "n", "19.11")
    ,@(353, "2025-06-03 23:35:51", 67, "n", "17.04")
    ,@(354, "2025-06-03 23:36:00", 72, "n", "18.46")
    ,@(355, "2025-06-03 23:36:06", 58, "n", "35.4")
    ,@(356, "2025-06-03 23:36:07", 53, "n", "40.68")
    ,@(357, "2025-06-03 23:56:35", 2, "s", "31.51")
    ,@(358, "2025-06-03 23:56:36", 4, "s", "45.77")
    ,@(359, "2025-06-03 23:56:37", 5, "s", "41.46")
    ,@(360, "2025-06-03 23:56:38", 8, "s", "18.28")
    ,@(361, "2025-06-03 23:56:58", 14, "s", "27.01")
    ,@(362, "2025-06-03 23:57:02", 15, "s", "33.04")
    ,@(363, "2025-06-03 23:57:05", 20, "s", "40.92")
    ,@(364, "2025-06-03 23:57:30", 36, "s", "31.51")
    ,@(365, "2025-06-03 23:57:31", 39, "s", "41.78")
    ,@(366, "2025-06-03 23:57:32", 40, "s", "18.28")
    ,@(367, "2025-06-03 23:57:35", 43, "s", "19.10")
    ,@(368, "2025-06-03 23:57:36", 45, "s", "42.31")
    ,@(369, "2025-06-03 23:57:37", 46, "s", "17.96")
    ,@(370, "2025-06-03 23:58:36", 49, "s", "27.01")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $time = $row[1]
    $vid = $row[2]
    $dtype = $row[3]
    $dval = $row[4]

    $ws.Cells.Item($r, 1).Value = $time
    $ws.Cells.Item($r, 2).Value = $vid
    $ws.Cells.Item($r, 3).Style = "Normal"  # materialize blank Class cell, no class recorded yet
    if ($dtype -eq "n") {
        $ws.Cells.Item($r, 4).Value = [double]$dval
    } else {
        # force text storage so the trailing-zero-preserving string round-trips verbatim
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $dval
        $ws.Cells.Item($r, 4).Style = "Normal"
    }
}

Write-Host "Added object detector selection rows 301-370; normalized D263:D300 to numeric."
